$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that sat right after the
#    "Staging URL => http://aidtheplanet.herokuapp.com" run.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Locate the "Tools" line in the skills section and widen it with the
#    extra technologies (React / Bootstrap / Ratchet CSS Framework).
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("jQuery, AngularJS, Backbone.js, Ember.js, Git", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the skills run to edit"
}

$base = $r.Start

# "jQuery" is the first 6 characters of the matched range.
$jqEnd = $base + 6

# Insert ", React" directly after "jQuery" as its own run (by nudging the
# character formatting of the freshly inserted text so Word has to give it
# a dedicated <w:r>/<w:rPr> instead of folding it into the neighboring run).
$insReact = $d.Range($jqEnd, $jqEnd)
$insReact.InsertAfter(", React")
$reactLen = ", React".Length
$reactRange = $d.Range($jqEnd, $jqEnd + $reactLen)
$reactRange.Font.Bold = 1
$reactRange.Font.Bold = 0

# "Ember.js" originally ended 40 characters after $base; the ", React"
# insertion above shifted everything after it forward by its own length.
$emberEnd = $base + 40 + $reactLen

$insFrameworks = $d.Range($emberEnd, $emberEnd)
$insFrameworks.InsertAfter(", Bootstrap, Ratchet CSS Framework")
$fwLen = ", Bootstrap, Ratchet CSS Framework".Length
$fwRange = $d.Range($emberEnd, $emberEnd + $fwLen)
$fwRange.Font.Bold = 1
$fwRange.Font.Bold = 0

# Put the "_GoBack" bookmark back, now sitting right before the remaining
# ", " run (immediately after the newly inserted frameworks text).
$bmPos = $emberEnd + $fwLen
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))
